$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127; this shifts the existing rows 127..199
# down to 128..200 (and bumps the sheet dimension to A1:R200 automatically).
$ws.Rows("127").Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Range("A127").Value = 4
$ws.Range("B127").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C127").Value = "Los Lagos"
$ws.Range("D127").Value = 45097
$ws.Range("E127").Value = 10
$ws.Range("F127").Value = 100112052
$ws.Range("G127").Value = "Albahaca"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 90
$ws.Range("K127").Value = 5000
$ws.Range("L127").Value = 5000
$ws.Range("M127").Value = 5000
$ws.Range("N127").Value = "$/paquete"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 5000
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
